# Trade #8 closed at 2026-02-16 22:52:44 - base_strategy UP +0.000%
# Appends the new trade row (row 9) to both the "All Trades" and
# "base_strategy" worksheets, matching the existing OPEN-trade row layout.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 9

    $ws.Cells.Item($row, 1).Value = 8

    # Column B holds a literal "YYYY-MM-DD" label (not a real date), so force
    # text formatting first to stop Excel's automatic date conversion, then
    # strip the resulting style back off so the cell stays on the default
    # (unstyled) format, same as every other row in the table.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "22:52:44"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 49.999998
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
